# "sửa file mẫu import data mkt" — update the Mau_nhap_data_mkt.xlsx template:
#  - remove the 6th/last sample data row (row 10) on Sheet1
#  - make Sheet1 the active/selected tab again (it had drifted to "Quy Tắc")
#  - update the remembered selections on both sheets

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Sheet1
$ws2 = $wb.Worksheets.Item(2)   # Quy Tắc

# Drop the extra sample row (B10:I10) — "Nguyễn Huy Hoàng 4" / 0355230188 / ...
# This also drops the now-unused shared strings (0355230188, Nguyễn Huy Hoàng 4,
# MayBach 66, Mua thêm nguyên vật liệu, 147 Cầu Giấy Hà Nội, 4/11/2019 08:20:50)
# automatically on save.
$ws1.Range("B10:I10").ClearContents()

# Sheet1 becomes the active sheet/tab again.
$ws1.Activate()

# Restore view state: Sheet1 keeps its B1 top-left scroll and its selection
# moves to D14; Quy Tắc keeps its existing G16 selection untouched.
$aw = $excel.ActiveWindow
$aw.ScrollColumn = 2
$aw.ScrollRow = 1
$ws1.Range("D14").Select()

$aw.Left = 3420
$aw.Top = 3420
